$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-looking decimal numbers as literal TEXT in the source data
# (e.g. "1.002", "0.4714"). Mark the range as Text before writing so COM does not
# auto-coerce these into numeric cells, which would change both type and the
# on-screen representation (trailing zeros, thousand-dot groupings like "30.560.96").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "30.560.96"
$ws.Cells.Item(2, 5).Value = "  -0.05%  "

$ws.Cells.Item(3, 4).Value = "1.886.31"
$ws.Cells.Item(3, 5).Value = "  +0.74%  "

$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$ws.Cells.Item(5, 4).Value = "244.81"
$ws.Cells.Item(5, 5).Value = "  -1.01%  "

$ws.Cells.Item(6, 4).Value = "1.002"
$ws.Cells.Item(6, 5).Value = "  +0.13%  "

$ws.Cells.Item(7, 4).Value = "0.4714"
$ws.Cells.Item(7, 5).Value = "  -0.51%  "

$ws.Cells.Item(8, 4).Value = "0.2911"
$ws.Cells.Item(8, 5).Value = "  +0.16%  "

$ws.Cells.Item(9, 4).Value = "0.06490"
$ws.Cells.Item(9, 5).Value = "  +0.23%  "

$ws.Cells.Item(10, 4).Value = "22.44"
$ws.Cells.Item(10, 5).Value = "  +1.37%  "

$ws.Cells.Item(11, 4).Value = "0.07758"
$ws.Cells.Item(11, 5).Value = "  +0.63%  "

$ws.Cells.Item(12, 4).Value = "1.888.13"
$ws.Cells.Item(12, 5).Value = "  +1.04%  "

$ws.Cells.Item(13, 4).Value = "0.7374"
$ws.Cells.Item(13, 5).Value = "  -0.61%  "

$ws.Cells.Item(14, 4).Value = "95.87"
$ws.Cells.Item(14, 5).Value = "  -0.68%  "

$ws.Cells.Item(15, 4).Value = "5.184"
$ws.Cells.Item(15, 5).Value = "  +0.67%  "

$ws.Cells.Item(16, 4).Value = "283.15"
$ws.Cells.Item(16, 5).Value = "  +3.55%  "

$ws.Cells.Item(17, 4).Value = "30.630.48"
$ws.Cells.Item(17, 5).Value = "  +0.23%  "

$ws.Cells.Item(18, 4).Value = "13.07"
$ws.Cells.Item(18, 5).Value = "  -2.06%  "

$ws.Cells.Item(19, 5).Value = "  +0.15%  "

$ws.Cells.Item(20, 4).Value = "0.000007499"
$ws.Cells.Item(20, 5).Value = "  +0.06%  "

$ws.Cells.Item(21, 4).Value = "2.131.33"
$ws.Cells.Item(21, 5).Value = "  +0.70%  "

$ws.Cells.Item(22, 5).Value = "  +0.13%  "

$ws.Cells.Item(23, 4).Value = "5.263"
$ws.Cells.Item(23, 5).Value = "  +0.16%  "

$ws.Cells.Item(24, 4).Value = "6.255"
$ws.Cells.Item(24, 5).Value = "  +1.07%  "

$ws.Cells.Item(25, 4).Value = "9.156"
$ws.Cells.Item(25, 5).Value = "  -0.83%  "

$ws.Cells.Item(26, 4).Value = "164.29"
$ws.Cells.Item(26, 5).Value = "  +0.59%  "

$ws.Cells.Item(27, 4).Value = "18.84"
$ws.Cells.Item(27, 5).Value = "  +0.28%  "

$ws.Cells.Item(28, 4).Value = "1.900"
$ws.Cells.Item(28, 5).Value = "  -0.83%  "

$ws.Cells.Item(29, 4).Value = "1.353"
$ws.Cells.Item(29, 5).Value = "  +0.44%  "

$ws.Cells.Item(30, 4).Value = "0.09731"
$ws.Cells.Item(30, 5).Value = "  -2.50%  "

$ws.Cells.Item(31, 4).Value = "1.476"
$ws.Cells.Item(31, 5).Value = "  -1.97%  "

$ws.Cells.Item(32, 4).Value = "4.295"
$ws.Cells.Item(32, 5).Value = "  +0.02%  "

$ws.Cells.Item(33, 4).Value = "4.132"
$ws.Cells.Item(33, 5).Value = "  +0.63%  "

$ws.Cells.Item(34, 4).Value = "0.04873"
$ws.Cells.Item(34, 5).Value = "  +1.87%  "

$ws.Cells.Item(35, 4).Value = "1.127"
$ws.Cells.Item(35, 5).Value = "  +0.72%  "

$ws.Cells.Item(36, 4).Value = "0.6933"
$ws.Cells.Item(36, 5).Value = "  -0.46%  "

$ws.Cells.Item(37, 4).Value = "2.708"
$ws.Cells.Item(37, 5).Value = "  -0.31%  "

$ws.Cells.Item(38, 4).Value = "0.01900"
$ws.Cells.Item(38, 5).Value = "  +2.77%  "

$ws.Cells.Item(39, 4).Value = "2.831"
$ws.Cells.Item(39, 5).Value = "  +2.91%  "

$ws.Cells.Item(40, 4).Value = "75.40"
$ws.Cells.Item(40, 5).Value = "  +2.67%  "

$ws.Cells.Item(41, 4).Value = "6.198"
$ws.Cells.Item(41, 5).Value = "  +0.13%  "

$ws.Cells.Item(42, 4).Value = "2.007"
$ws.Cells.Item(42, 5).Value = "  +2.00%  "

$ws.Cells.Item(43, 4).Value = "0.4264"
$ws.Cells.Item(43, 5).Value = "  +2.05%  "

$ws.Cells.Item(44, 4).Value = "1.002"
$ws.Cells.Item(44, 5).Value = "  +0.11%  "

$ws.Cells.Item(45, 4).Value = "0.8231"
$ws.Cells.Item(45, 5).Value = "  -1.19%  "

$ws.Cells.Item(46, 4).Value = "101.29"
$ws.Cells.Item(46, 5).Value = "  -1.34%  "

$ws.Cells.Item(47, 4).Value = "9.524"
$ws.Cells.Item(47, 5).Value = "  +2.41%  "

$ws.Cells.Item(48, 4).Value = "35.40"
$ws.Cells.Item(48, 5).Value = "  +0.07%  "

$ws.Cells.Item(49, 4).Value = "6.976"
$ws.Cells.Item(49, 5).Value = "  +0.23%  "

$ws.Cells.Item(50, 4).Value = "912.63"
$ws.Cells.Item(50, 5).Value = "  -1.46%  "

$ws.Cells.Item(51, 4).Value = "0.05751"
$ws.Cells.Item(51, 5).Value = "  +1.84%  "
